$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.612.59"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").Value = "'1.585.55"
$ws.Range("E3").Value = "  -0.38%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'207.86"
$ws.Range("E5").Value = "  +0.41%  "
$ws.Range("D6").Value = "'0.499"
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D8").Value = "'22.23"
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("D9").Value = "'0.250"
$ws.Range("E9").Value = "  -0.35%  "
$ws.Range("D11").Value = "'0.0867"
$ws.Range("E11").Value = "  -0.58%  "
$ws.Range("D12").Value = "'1.811.72"
$ws.Range("E12").Value = "  -0.41%  "
$ws.Range("D13").Value = "'1.604.37"
$ws.Range("E13").Value = "  +0.62%  "
$ws.Range("E14").Value = "  -0.52%  "
$ws.Range("D15").Value = "'0.526"
$ws.Range("E15").Value = "  -1.90%  "
$ws.Range("D16").Value = "'27.622.42"
$ws.Range("E16").Value = "  +0.43%  "
$ws.Range("D17").Value = "'63.01"
$ws.Range("E17").Value = "  -0.25%  "
$ws.Range("D18").Value = "'216.95"
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("E19").Value = "  +0.27%  "
$ws.Range("D20").Value = "'7.30"
$ws.Range("E20").Value = "  -0.67%  "
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("D22").Value = "'4.13"
$ws.Range("E22").Value = "  -1.35%  "
$ws.Range("D23").Value = "'9.72"
$ws.Range("E23").Value = "  +0.65%  "
$ws.Range("D24").Value = "'1.98"
$ws.Range("E24").Value = "  -0.86%  "
$ws.Range("D25").Value = "'153.52"
$ws.Range("E25").Value = "  -0.91%  "
$ws.Range("D26").Value = "'7.03"
$ws.Range("E26").Value = "  +5.05%  "
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("D28").Value = "'15.04"
$ws.Range("E28").Value = "  +0.35%  "
$ws.Range("E29").Value = "  -0.73%  "
$ws.Range("E30").Value = "  -0.07%  "
$ws.Range("D31").Value = "'0.0472"
$ws.Range("E31").Value = "  +1.03%  "
$ws.Range("E32").Value = "  -2.06%  "
$ws.Range("D33").Value = "'1.371.76"
$ws.Range("E33").Value = "  +1.48%  "
$ws.Range("E34").Value = "  +0.57%  "
$ws.Range("D35").Value = "'1.54"
$ws.Range("E35").Value = "  +0.09%  "
$ws.Range("D36").Value = "'0.968"
$ws.Range("E36").Value = "  +1.47%  "
$ws.Range("E37").Value = "  +0.11%  "
$ws.Range("E38").Value = "  +2.40%  "
$ws.Range("D39").Value = "'0.532"
$ws.Range("E39").Value = "  -0.78%  "
$ws.Range("E40").Value = "  +1.71%  "
$ws.Range("E41").Value = "  +0.10%  "
$ws.Range("D42").Value = "'0.973"
$ws.Range("E42").Value = "  +1.54%  "
$ws.Range("D43").Value = "'64.16"
$ws.Range("E43").Value = "  +0.68%  "
$ws.Range("E44").Value = "  +4.17%  "
$ws.Range("E45").Value = "  +0.98%  "
$ws.Range("E46").Value = "  -1.67%  "
$ws.Range("D47").Value = "'1.723.81"
$ws.Range("E47").Value = "  -0.23%  "
$ws.Range("D48").Value = "'85.96"
$ws.Range("E48").Value = "  -1.96%  "
$ws.Range("D49").Value = "'0.0₆0101"
$ws.Range("E49").Value = "  +2.35%  "
$ws.Range("D50").Value = "'0.0961"
$ws.Range("E50").Value = "  -0.46%  "
$ws.Range("D51").Value = "'0.0494"
$ws.Range("E51").Value = "  -0.54%  "
